# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---- Sheet "Home win" ----
$ws = $wb.Worksheets.Item("Home win")
$ws.Range("A2").Value = "18-12-2024 00:00"
$ws.Range("B2").Value = "BOLIVIA"
$ws.Range("C2").Value = "PRIMERA DIVISIÓN"
$ws.Range("D2").Value = "Guabirá - Royal Pari"
$ws.Range("E2").Value = 73.3
$ws.Range("F2").Value = 1.75
$ws.Rows.Item(3).Delete()

# ---- Sheet "Draw" ----
$ws = $wb.Worksheets.Item("Draw")
$ws.Range("A3").Value = "18-12-2024 10:30"
$ws.Range("B3").Value = "WORLD"
$ws.Range("C3").Value = "AFF CHAMPIONSHIP"
$ws.Range("D3").Value = "Myanmar - Laos"
$ws.Range("E3").Value = 60

# ---- Sheet "Btts" ----
$ws = $wb.Worksheets.Item("Btts")
$ws.Range("A2").Value = "17-12-2024 19:45"
$ws.Range("B2").Value = "ENGLAND"
$ws.Range("C2").Value = "NATIONAL LEAGUE"
$ws.Range("D2").Value = "Rochdale - Tamworth"
$ws.Range("E2").Value = 84
$ws.Range("F2").Value = 1.8

$ws.Range("A3").Value = "17-12-2024 20:00"
$ws.Range("B3").Value = "NETHERLANDS"
$ws.Range("C3").Value = "KNVB BEKER"
$ws.Range("D3").Value = "MVV - Feyenoord"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 2

$ws.Range("A4").Value = "17-12-2024 08:30"
$ws.Range("B4").Value = "INDONESIA"
$ws.Range("C4").Value = "LIGA 1"
$ws.Range("D4").Value = "Malut United - PSM Makassar"
$ws.Range("E4").Value = 76
$ws.Range("F4").Value = 1.85

$ws.Range("A5").Value = "17-12-2024 12:30"
$ws.Range("B5").Value = "WORLD"
$ws.Range("C5").Value = "AFF CHAMPIONSHIP"
$ws.Range("D5").Value = "Singapore - Thailand"
$ws.Range("E5").Value = 86.7
$ws.Range("F5").Value = 1.95

$ws.Range("A6").Value = "18-12-2024 17:30"
$ws.Range("B6").Value = "ITALY"
$ws.Range("C6").Value = "COPPA ITALIA"
$ws.Range("D6").Value = "Atalanta - Cesena"
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 2.1

$ws.Range("A7").Value = "18-12-2024 19:00"
$ws.Range("B7").Value = "NETHERLANDS"
$ws.Range("C7").Value = "KNVB BEKER"
$ws.Range("D7").Value = "ASWH - Heerenveen"
$ws.Range("E7").Value = 76
$ws.Range("F7").Value = 1.7

$ws.Range("A8").Value = "18-12-2024 19:00"
$ws.Range("B8").Value = "NETHERLANDS"
$ws.Range("C8").Value = "KNVB BEKER"
$ws.Range("D8").Value = "AFC Amsterdam - Utrecht"
$ws.Range("E8").Value = 76
$ws.Range("F8").Value = 1.91

$ws.Range("A9").Value = "18-12-2024 15:00"
$ws.Range("B9").Value = "ROMANIA"
$ws.Range("C9").Value = "CUPA ROMÂNIEI"
$ws.Range("D9").Value = "Politehnica Iasi - AFC Hermannstadt"
$ws.Range("E9").Value = 80
$ws.Range("F9").Value = 1.91

$ws.Range("A10").Value = "18-12-2024 18:00"
$ws.Range("B10").Value = "SPAIN"
$ws.Range("C10").Value = "SEGUNDA DIVISIÓN"
$ws.Range("D10").Value = "Racing Ferrol - Almeria"
$ws.Range("E10").Value = 76
$ws.Range("F10").Value = 1.77

$ws.Range("A11").Value = "24-11-2024 16:00"
$ws.Range("B11").Value = "SPAIN"
$ws.Range("C11").Value = "SEGUNDA DIVISIÓN RFEF - GROUP 1"
$ws.Range("D11").Value = "Bergantiños - Deportivo La Coruña II"
$ws.Range("E11").Value = 80
$ws.Range("F11").Value = 1.8

# ---- Sheet "Over_Under" ----
$ws = $wb.Worksheets.Item("Over_Under")
$ws.Range("A2").Value = "17-12-2024 19:00"
$ws.Range("B2").Value = "NETHERLANDS"
$ws.Range("C2").Value = "KNVB BEKER"
$ws.Range("D2").Value = "Waalwijk - Cambuur"
$ws.Range("E2").Value = 70
$ws.Range("F2").Value = 1.7
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 2.75

$ws.Range("A3").Value = "17-12-2024 19:00"
$ws.Range("B3").Value = "BOLIVIA"
$ws.Range("C3").Value = "PRIMERA DIVISIÓN"
$ws.Range("D3").Value = "Gualberto Villarroel SJ - Blooming"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 1.73
$ws.Range("G3").Value = 66.7
$ws.Range("H3").Value = 2.75

$ws.Range("A6").Value = "18-12-2024 00:00"
$ws.Range("B6").Value = "BOLIVIA"
$ws.Range("C6").Value = "PRIMERA DIVISIÓN"
$ws.Range("D6").Value = "Guabirá - Royal Pari"
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 1.77
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = 2.8
